$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 273; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $old = $cell.Value2
    $new = [Math]::Round($old * 7, 2)
    $cell.Value2 = $new
}
